$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 5977.3335
$ws.Range("I76").Value = 4928
$ws.Range("J76").Value = 9650
$ws.Range("K76").Value = 4928
$ws.Range("L76").Value = 9650
$ws.Range("M76").Value = -4613
$ws.Range("N76").Value = -10280

# Row 79
$ws.Range("H79").Value = 5977.3335
$ws.Range("I79").Value = 4928
$ws.Range("J79").Value = 9650
$ws.Range("K79").Value = 4928
$ws.Range("L79").Value = 9650
$ws.Range("M79").Value = -3836
$ws.Range("N79").Value = -11834

# Row 138
$ws.Range("H138").Value = 6829.0864
$ws.Range("I138").Value = 3155.3845
$ws.Range("J138").Value = 7890.378
$ws.Range("K138").Value = 9466.1535
$ws.Range("L138").Value = 23671.134
$ws.Range("M138").Value = -4326.1535
$ws.Range("N138").Value = -33951.134

$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 57250
$ws.Range("I37").Value = 49500
$ws.Range("K37").Value = 49500
$ws.Range("M37").Value = -49227

# Row 88
$ws.Range("H88").Value = 1902
$ws.Range("I88").Value = 1781.3334
$ws.Range("J88").Value = 1992.5
$ws.Range("K88").Value = 1781.3334
$ws.Range("L88").Value = 1992.5
$ws.Range("M88").Value = -1375.3334
$ws.Range("N88").Value = -2804.5

# Row 91
$ws.Range("H91").Value = 1902
$ws.Range("I91").Value = 1781.3334
$ws.Range("J91").Value = 1992.5
$ws.Range("K91").Value = 1781.3334
$ws.Range("L91").Value = 1992.5
$ws.Range("M91").Value = -377.3334
$ws.Range("N91").Value = -4800.5

# Row 97
$ws.Range("H97").Value = 1292.8572
$ws.Range("I97").Value = 571.7273
$ws.Range("J97").Value = 3937
$ws.Range("K97").Value = 571.7273
$ws.Range("L97").Value = 3937
$ws.Range("M97").Value = -75.72730000000001
$ws.Range("N97").Value = -4929

# Row 102
$ws.Range("H102").Value = 3131.524
$ws.Range("I102").Value = 2934.8948
$ws.Range("K102").Value = 2934.8948
$ws.Range("M102").Value = -1312.8948

# Row 132
$ws.Range("H132").Value = 62504610
$ws.Range("I132").Value = 4454.154
$ws.Range("K132").Value = 13362.462
$ws.Range("M132").Value = -10832.462

$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Range("H10").Value = 3436.25
$ws.Range("J10").Value = 350
$ws.Range("L10").Value = 350
$ws.Range("N10").Value = -630

# Row 25
$ws.Range("H25").Value = 2353.875
$ws.Range("I25").Value = 2589.2
$ws.Range("J25").Value = 1961.6666
$ws.Range("K25").Value = 2589.2
$ws.Range("L25").Value = 1961.6666
$ws.Range("M25").Value = -2354.2
$ws.Range("N25").Value = -2431.6666

# Row 54
$ws.Range("H54").Value = 2250
$ws.Range("I54").Value = 2250
$ws.Range("K54").Value = 2250
$ws.Range("M54").Value = -1766

# Row 105
$ws.Range("H105").Value = 1831.3334
$ws.Range("I105").Value = 1500
$ws.Range("J105").Value = 1997
$ws.Range("K105").Value = 1500
$ws.Range("L105").Value = 1997
$ws.Range("M105").Value = 247
$ws.Range("N105").Value = -5491

# Row 134
$ws.Range("H134").Value = 1560.4231
$ws.Range("I134").Value = 1560.4231
$ws.Range("K134").Value = 4681.2693
$ws.Range("M134").Value = -2146.2693

$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 849.6
$ws.Range("I13").Value = 999
$ws.Range("K13").Value = 999
$ws.Range("M13").Value = -860

# Row 16
$ws.Range("H16").Value = 2013.7646
$ws.Range("I16").Value = 2053.6924
$ws.Range("J16").Value = 1884
$ws.Range("K16").Value = 2053.6924
$ws.Range("L16").Value = 1884
$ws.Range("M16").Value = -1766.6924
$ws.Range("N16").Value = -2458

# Row 17
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# Row 105
$ws.Range("H105").Value = 12755.6
$ws.Range("I105").Value = 3060.5557
$ws.Range("K105").Value = 3060.5557
$ws.Range("M105").Value = -1313.5557

# Row 113
$ws.Range("H113").Value = 2013.7646
$ws.Range("I113").Value = 2053.6924
$ws.Range("J113").Value = 1884
$ws.Range("K113").Value = 2053.6924
$ws.Range("L113").Value = 1884
$ws.Range("M113").Value = 116.3076000000001
$ws.Range("N113").Value = -6224

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1616.0513
$ws.Range("I5").Value = 1029.9642
$ws.Range("K5").Value = 3089.8926
$ws.Range("M5").Value = -2977.8926

# Row 32
$ws.Range("H32").Value = 2123.25
$ws.Range("J32").Value = 2664.3333
$ws.Range("L32").Value = 7992.999899999999
$ws.Range("N32").Value = -8558.999899999999

# Row 131
$ws.Range("H131").Value = 1755.8334
$ws.Range("J131").Value = 1850.4166
$ws.Range("L131").Value = 5551.2498
$ws.Range("N131").Value = -15631.2498

# Row 135
$ws.Range("H135").Value = 1616.0513
$ws.Range("I135").Value = 1029.9642
$ws.Range("K135").Value = 9269.677799999999
$ws.Range("M135").Value = -6734.677799999999

$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 31499.5
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

# Row 122
$ws.Range("H122").Value = 3070.25
$ws.Range("I122").Value = 2891
$ws.Range("K122").Value = 8673
$ws.Range("M122").Value = -6223

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4623.933
$ws.Range("I7").Value = 4247.3335
$ws.Range("K7").Value = 4247.3335
$ws.Range("M7").Value = -4135.3335

# Row 46
$ws.Range("H46").Value = 1532
$ws.Range("I46").Value = 899.3333
$ws.Range("K46").Value = 899.3333
$ws.Range("M46").Value = -711.3333

# Row 82
$ws.Range("H82").Value = 2861.6875
$ws.Range("I82").Value = 1072
$ws.Range("K82").Value = 1072
$ws.Range("M82").Value = -711

# Row 85
$ws.Range("H85").Value = 2861.6875
$ws.Range("I85").Value = 1072
$ws.Range("K85").Value = 1072
$ws.Range("M85").Value = 176

# Row 93
$ws.Range("H93").Value = 1587.6666
$ws.Range("I93").Value = 1342.4375
$ws.Range("J93").Value = 3549.5
$ws.Range("K93").Value = 1342.4375
$ws.Range("L93").Value = 3549.5
$ws.Range("M93").Value = -94.4375
$ws.Range("N93").Value = -6045.5

# Row 126
$ws.Range("H126").Value = 4623.933
$ws.Range("I126").Value = 4247.3335
$ws.Range("K126").Value = 12742.0005
$ws.Range("M126").Value = -10272.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 6000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# Row 33
$ws.Range("H33").Value = 8173.5713
$ws.Range("J33").Value = 13999
$ws.Range("L33").Value = 13999
$ws.Range("N33").Value = -14499

# Row 36
$ws.Range("H36").Value = 8173.5713
$ws.Range("J36").Value = 13999
$ws.Range("L36").Value = 13999
$ws.Range("N36").Value = -14499

# Row 52
$ws.Range("H52").Value = 27498.715
$ws.Range("I52").Value = 24749.25
$ws.Range("K52").Value = 24749.25
$ws.Range("M52").Value = -24523.25

# Row 76
$ws.Range("H76").Value = 77500
$ws.Range("I76").Value = 80000
$ws.Range("J76").Value = 75000
$ws.Range("K76").Value = 80000
$ws.Range("L76").Value = 75000
$ws.Range("M76").Value = -79685
$ws.Range("N76").Value = -75630

# Row 79
$ws.Range("H79").Value = 77500
$ws.Range("I79").Value = 80000
$ws.Range("J79").Value = 75000
$ws.Range("K79").Value = 80000
$ws.Range("L79").Value = 75000
$ws.Range("M79").Value = -78908
$ws.Range("N79").Value = -77184

# Row 96
$ws.Range("H96").Value = 8216.286
$ws.Range("I96").Value = 7854
$ws.Range("K96").Value = 7854
$ws.Range("M96").Value = -6481

# Row 107
$ws.Range("H107").Value = 749.5
$ws.Range("I107").Value = 688.2222
$ws.Range("J107").Value = 933.3333
$ws.Range("K107").Value = 2064.6666
$ws.Range("L107").Value = 2799.9999
$ws.Range("M107").Value = -144.6666
$ws.Range("N107").Value = -6639.9999

# Row 132
$ws.Range("H132").Value = 6433.628
$ws.Range("I132").Value = 6372.3335
$ws.Range("J132").Value = 6748.857
$ws.Range("K132").Value = 19117.0005
$ws.Range("L132").Value = 20246.571
$ws.Range("M132").Value = -16587.0005
$ws.Range("N132").Value = -25306.571

# Row 135
$ws.Range("H135").Value = 16722370
$ws.Range("J135").Value = 16722370
$ws.Range("L135").Value = 16722370
$ws.Range("N135").Value = -16732510
